# Removed radiobutton from input panel config.
# On the "Translation" sheet, two pairs of rows (each holding the
# "Manual"/"Defined" radiobutton option labels for an input panel) are
# removed. Removing the rows shifts all subsequent rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Delete rows 60 and 61 first (higher row numbers) so that the row numbers
# for the first pair (26 and 27) remain valid afterwards.
$ws.Rows.Item(60).Resize(2).Delete() | Out-Null

# Delete the first pair of rows (26 and 27).
$ws.Rows.Item(26).Resize(2).Delete() | Out-Null
